# "Compil - corrections staff et locaux"
# Add the new "Material deposit - Neutral support" row to the ADMIN sheet
# (row 21), matching the layout of the other location rows above it, and
# update the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADMIN")

# --- populate the new row's content -------------------------------------
$ws.Range("A21").Value = "Material deposit - Neutral support"
$ws.Range("B21").Value = "Dépôt de matérie -  Dépannage neutre"
$ws.Range("C21").Value = "TBD"
$ws.Range("D21").Formula = "=C21"

# --- copy formatting from comparable existing rows -----------------------
# A21/B21 take the same look as A20/B20 (left column style + label style).
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C21/D21 (the "TBD" location + its mirrored formula) match row 17, the
# other "TBD" location row.
$ws.Range("C17:D17").Copy()
$ws.Range("C21:D21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row grows to fit the two-line label, same as the other wrapped rows.
$ws.Rows.Item(21).RowHeight = 34

# --- update the view / selection to match ---------------------------------
$ws.Range("D21").Select()
$excel.ActiveWindow.ScrollColumn = 2
